# Actualización automática 2025-11-27 15:30:08
#
# The workbook stores plain cached values (no formulas), so the monthly
# "PORCELANATO" (M4) and "LAVABOS" (I31) sales figures for asesor
# GUERRERO FAREZ FABIAN MAURICIO are corrected on the detail sheet, and the
# matching downstream roll-ups on the "VENTA MENSUAL" and
# "CUMPLIMIENTO MENSUAL" sheets are updated to stay consistent.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO ------------------------------------------------------
# AGUILAR REYES CESAR VINICIO - PORCELANATO (noviembre)
$wsGrupo.Range("M4").Value = 8294.530000000001
# MATUTE GUANOLIQUE DOLORES MATILDE - LAVABOS (noviembre)
$wsGrupo.Range("I31").Value = 79.2
# PEREZ ROSALES EDGAR RICARDO - PORCELANATO (noviembre)
$wsGrupo.Range("M40").Value = 0
# Summary row: count of non-zero PORCELANATO entries out of 54
$wsGrupo.Range("M56").Value = "16 de 54"

# --- VENTA MENSUAL ----------------------------------------------------------
# AGUILAR REYES CESAR VINICIO - noviembre total
$wsMensual.Range("F4").Value = 11015.52
# MATUTE GUANOLIQUE DOLORES MATILDE - noviembre total
$wsMensual.Range("F31").Value = 79.2
# PEREZ ROSALES EDGAR RICARDO - noviembre total
$wsMensual.Range("F40").Value = 0
# Column total for noviembre
$wsMensual.Range("F60").Value = 80203.93000000001

# --- CUMPLIMIENTO MENSUAL ----------------------------------------------------
# LAVABOS row
$wsCumplimiento.Range("D7").Value = 702
$wsCumplimiento.Range("E7").Value = 618
$wsCumplimiento.Range("F7").Value = 0.5318181818181819

# PORCELANATO row
$wsCumplimiento.Range("D12").Value = 41166.81
$wsCumplimiento.Range("E12").Value = 23777.19
$wsCumplimiento.Range("F12").Value = 0.6338816518847006

# TOTAL row
$wsCumplimiento.Range("D14").Value = 77221.22
$wsCumplimiento.Range("E14").Value = 21735.03685923838
$wsCumplimiento.Range("F14").Value = 0.7803571239547221
